# "Data Munging & Plot Work" — mark progress on the first chart task
# (Total Carbon Offset by State bar plot) in the "Charts to Make" sheet.

$wb = $excel.ActiveWorkbook

$wsCharts = $wb.Worksheets.Item("Charts to Make")
$wsMeta   = $wb.Worksheets.Item("Metadata")

# Row 5 = "Total Carbon Offset by State" (Bar chart) task row.
# "Completed?" (col E) moves from "No" to "In Progress" -> highlight yellow.
$eCell = $wsCharts.Range("E5")
$eCell.Value = "In Progress"
$eCell.Interior.Color = 65535        # RGB(255,255,0) yellow

# "Completed?3" (col I) moves from "No" to "Yes" -> highlight green,
# matching the existing "Completed?2" (col G) styling for this row.
$iCell = $wsCharts.Range("I5")
$iCell.Value = "Yes"
$iCell.Interior.Color = 5287936      # RGB(0,176,80) green

# Leave the cursor/selection where the author last left it when saving.
$wsMeta.Range("B9").Select()
$wsCharts.Activate()
$wsCharts.Range("F19").Select()
